$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1206
$ws1.Range("F5").Value = 839
$ws1.Range("F7").Value = 1579
$ws1.Range("F9").Value = 1078
$ws1.Range("F12").Value = 212
$ws1.Range("F13").Value = 65
$ws1.Range("F14").Value = 543
$ws1.Range("F19").Value = 304
$ws1.Range("F24").Value = 795

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 207
$ws2.Range("F7").Value = 75
$ws2.Range("F8").Value = 602

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1206
$ws4.Range("F7").Value = 839
$ws4.Range("F9").Value = 1579
$ws4.Range("F11").Value = 1078
$ws4.Range("F14").Value = 212
$ws4.Range("F15").Value = 65
$ws4.Range("F16").Value = 543
$ws4.Range("F23").Value = 304
$ws4.Range("F25").Value = 207
$ws4.Range("F26").Value = 207
$ws4.Range("F31").Value = 795
$ws4.Range("F33").Value = 75
$ws4.Range("F35").Value = 602
